$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Source: NBS" block (rows 21-22) down to rows 27-28 ---
# Capture existing values first (they are the same text, just relocated)
$nbsVal = $ws.Range("A21").Value2
$nbsDescVal = $ws.Range("A22").Value2

$ws.Range("A27").Value2 = $nbsVal
$ws.Range("A28").Value2 = $nbsDescVal

# Clear the old cells (rows 21-22 are no longer used)
$ws.Range("A21").Value2 = ""
$ws.Range("A22").Value2 = ""
$ws.Range("A21").ClearContents()
$ws.Range("A22").ClearContents()

# --- New size-classification table (rows 16-20) ---
$ws.Range("B16").Value2 = "Number of employees"
$ws.Range("C16").Value2 = "Assets (local currency, unless noted otherwise)"
$ws.Range("D16").Value2 = "Turnover (local currency, unless noted otherwise)"

$ws.Range("A17").Value2 = "Micro"
$ws.Range("B17").Value2 = "1-4"
$ws.Range("C17").Value2 = ""
$ws.Range("D17").Value2 = ""

$ws.Range("A18").Value2 = "Small"
$ws.Range("B18").Value2 = "5-9"
$ws.Range("C18").Value2 = ""
$ws.Range("D18").Value2 = ""

$ws.Range("A19").Value2 = "Medium"
$ws.Range("B19").Value2 = "10-49"
$ws.Range("C19").Value2 = ""
$ws.Range("D19").Value2 = ""

$ws.Range("A20").Value2 = "Large"
$ws.Range("B20").Value2 = ">49"
$ws.Range("C20").Value2 = ""
$ws.Range("D20").Value2 = ""
